# Add two new rows (5 and 6) of event data to the "my_Events" sheet,
# extending the used range from A1:F4 to A1:F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: new event "Adrianas Bday" -------------------------------------
# All values in this table are stored as plain text (not numbers), even
# though some look numeric (IDs, day/month/reminder numbers). Force text
# by switching the cell format to Text ("@") before writing the value,
# then reset the visible format back to Normal/General so no stray
# number-format override is left behind on the cell.
$row5 = $ws.Range("A5:F5")
$row5.NumberFormat = "@"
$ws.Range("A5").Value = "3667518741744194029"
$ws.Range("B5").Value = "13"
$ws.Range("C5").Value = "11"
$ws.Range("D5").Value = "2022"
$ws.Range("E5").Value = "Adrianas Bday"
$ws.Range("F5").Value = "14"
$row5.Style = "Normal"

# --- Row 6: a fresh blank template row (only the generated ID is filled) --
$row6 = $ws.Range("A6:F6")
$row6.NumberFormat = "@"
$ws.Range("A6").Value = "1510441034453422573"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$row6.Style = "Normal"
